$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

$ws.Range("D20:D21").Select() | Out-Null
